$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with same style as the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate time_taken values for each data row (2-72)
$timeTaken = @{
    2 = "2021-10-05 13:41:35.633165"
    3 = "2021-10-05 13:41:35.633177"
    4 = "2021-10-05 13:41:35.633180"
    5 = "2021-10-05 13:41:35.633183"
    6 = "2021-10-05 13:41:35.633186"
    7 = "2021-10-05 13:41:35.633189"
    8 = "2021-10-05 13:41:35.633191"
    9 = "2021-10-05 13:41:35.633194"
    10 = "2021-10-05 13:41:35.633196"
    11 = "2021-10-05 13:41:35.633199"
    12 = "2021-10-05 13:41:35.633202"
    13 = "2021-10-05 13:41:35.633204"
    14 = "2021-10-05 13:41:35.633206"
    15 = "2021-10-05 13:41:35.633209"
    16 = "2021-10-05 13:41:35.633211"
    17 = "2021-10-05 13:41:35.633214"
    18 = "2021-10-05 13:41:35.633217"
    19 = "2021-10-05 13:41:35.633219"
    20 = "2021-10-05 13:41:35.633222"
    21 = "2021-10-05 13:41:35.633224"
    22 = "2021-10-05 13:41:35.633227"
    23 = "2021-10-05 13:41:35.633229"
    24 = "2021-10-05 13:41:35.633232"
    25 = "2021-10-05 13:41:35.633234"
    26 = "2021-10-05 13:41:35.633237"
    27 = "2021-10-05 13:41:35.633239"
    28 = "2021-10-05 13:41:35.633242"
    29 = "2021-10-05 13:41:35.633244"
    30 = "2021-10-05 13:41:35.633247"
    31 = "2021-10-05 13:41:35.633249"
    32 = "2021-10-05 13:41:35.633252"
    33 = "2021-10-05 13:41:35.633254"
    34 = "2021-10-05 13:41:35.633257"
    35 = "2021-10-05 13:41:35.633260"
    36 = "2021-10-05 13:41:35.633262"
    37 = "2021-10-05 13:41:35.633265"
    38 = "2021-10-05 13:41:35.633267"
    39 = "2021-10-05 13:41:35.633270"
    40 = "2021-10-05 13:41:35.633272"
    41 = "2021-10-05 13:41:35.633275"
    42 = "2021-10-05 13:41:35.633277"
    43 = "2021-10-05 13:41:35.633280"
    44 = "2021-10-05 13:41:35.633282"
    45 = "2021-10-05 13:41:35.633285"
    46 = "2021-10-05 13:41:35.633287"
    47 = "2021-10-05 13:41:35.633290"
    48 = "2021-10-05 13:41:35.633292"
    49 = "2021-10-05 13:41:35.633295"
    50 = "2021-10-05 13:41:35.633297"
    51 = "2021-10-05 13:41:35.633300"
    52 = "2021-10-05 13:41:35.633302"
    53 = "2021-10-05 13:41:35.633304"
    54 = "2021-10-05 13:41:35.633307"
    55 = "2021-10-05 13:41:35.633310"
    56 = "2021-10-05 13:41:35.633312"
    57 = "2021-10-05 13:41:35.633315"
    58 = "2021-10-05 13:41:35.633317"
    59 = "2021-10-05 13:41:35.633320"
    60 = "2021-10-05 13:41:35.633322"
    61 = "2021-10-05 13:41:35.633324"
    62 = "2021-10-05 13:41:35.633327"
    63 = "2021-10-05 13:41:35.633329"
    64 = "2021-10-05 13:41:35.633332"
    65 = "2021-10-05 13:41:35.633334"
    66 = "2021-10-05 13:41:35.633338"
    67 = "2021-10-05 13:41:35.633341"
    68 = "2021-10-05 13:41:35.633343"
    69 = "2021-10-05 13:41:35.633345"
    70 = "2021-10-05 13:41:35.633348"
    71 = "2021-10-05 13:41:35.633351"
    72 = "2021-10-05 13:41:35.633353"
}

foreach ($row in 2..72) {
    $ws.Cells.Item($row, 6).Value = $timeTaken[$row]
}
